$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# Mapping of the comment text that must end up at each row AFTER the insert,
# built from the ORIGINAL (pre-insert) comment text one row above (comments do
# not physically move when a row is inserted, so we shift their text manually,
# working on the ORIGINAL snapshot so order of application does not matter).
$commentShift = @(
    @{ Row = 19; Text = "Name of the virtual machine on the Hypervisor/vCenter Server. [default: (Hostname)]" },
    @{ Row = 20; Text = "Name of the system if other than hostname [default: (Hostname)]" },
    @{ Row = 21; Text = "Unique username used to identify this VSC in its XMPP connection with VSD [default: vsc1]" },
    @{ Row = 23; Text = "Type of hypervisor environment where VMs will be instantiated. Use 'none' when skipping predeploy." },
    @{ Row = 24; Text = "Hostname or IP address of the hypervisor where VM  will be instantiated. In the case of deployment in a vCenter environment, this will be the FQDN of the vCenter Server" },
    @{ Row = 26; Text = "Network Bridge used for the management interface of a component or the BOF interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. This field can be overridden by defining the management network bridge separately in the component configuration. Defaults to the global setting [default: (global Bridge interface)]" },
    @{ Row = 27; Text = "Network Bridge used for the data path of a component or the Control interface on VSC. This will be a Distributed Virtual PortGroup (DVPG) when deploying on vCenter or a Linux network bridge when deploying on KVM. [default: (global Bridge interface)]" },
    @{ Row = 28; Text = "FQDN of the VSD or VSD cluster for this VSC" },
    @{ Row = 29; Text = "Private Management IP Address of VSC instances" },
    @{ Row = 30; Text = "Private Control IP Address of VSC Instances" },
    @{ Row = 31; Text = "Private Data Gateway IP Address of VSC Instances" },
    @{ Row = 32; Text = "List of route reflector IP addresses if present (List items separated by comma.)" },
    @{ Row = 34; Text = "Name of the vCenter Datacenter on which the VSC VM will be deployed. Defaults to the common vCenter Datacenter Name if not defined here. [default: (global vCenter Datacenter Name)]" },
    @{ Row = 35; Text = "Name of the vCenter Cluster on which the VSC VM will be deployed. Defaults to the common vCenter Cluster Name if not defined here. [default: (global vCenter Cluster Name)]" },
    @{ Row = 36; Text = "Requires ovftool 4.3. Reference to the host on the vCenter cluster on which to deploy Nuage components [default: (global vCenter Host Reference)]" },
    @{ Row = 37; Text = "Name of the vCenter Datastore on which the VSC VM will be deployed. Defaults to the common vCenter Datastore Name if not defined here. [default: (global vCenter Datastore Name)]" },
    @{ Row = 38; Text = "Optional path to a folder defined on vCenter where VM will be instantiated [default: (global vCenter VM folder)]" },
    @{ Row = 39; Text = "Optional path to a hosts and clusters folder defined on vCenter where VM will be instantiated" },
    @{ Row = 41; Text = "Name of image installed on OpenStack for VSC" },
    @{ Row = 42; Text = "Name of instance flavor installed on OpenStack for VSC" },
    @{ Row = 43; Text = "Name of availability zone on OpenStack for VSC" },
    @{ Row = 44; Text = "Name of management network on OpenStack for VSC" },
    @{ Row = 45; Text = "Name of management subnet on OpenStack for VSC" },
    @{ Row = 46; Text = "Name for Mgmt interface" },
    @{ Row = 47; Text = "Set of security groups to associate with Mgmt interface (List items separated by comma.)" },
    @{ Row = 48; Text = "Name of control network on OpenStack for VSC" },
    @{ Row = 49; Text = "Name of control subnet on OpenStack for VSC" },
    @{ Row = 50; Text = "Name for Control interface" },
    @{ Row = 51; Text = "Set of security groups to associate with Control interface (List items separated by comma.)" },
    @{ Row = 52; Text = "Name for Mgmt interface" },
    @{ Row = 53; Text = "Set of security groups to associate with Mgmt interface (List items separated by comma.)" },
    @{ Row = 55; Text = "Used in postdeploy and health workflows as expected values if non-zero [default: 0]" },
    @{ Row = 56; Text = "Used in postdeploy and health workflows as expected values if non-zero [default: 0]" },
    @{ Row = 57; Text = "Used in postdeploy and health workflows as expected values if non-zero [default: 0]" },
    @{ Row = 58; Text = "Used in postdeploy and health workflows as expected values if non-zero [default: 0]" },
    @{ Row = 59; Text = "Used in postdeploy and health workflows as expected values if non-zero [default: 0]" },
    @{ Row = 61; Text = "Ejabberd user id used to create the certificate" },
    @{ Row = 62; Text = "Path to VSC certificate key pem file" },
    @{ Row = 63; Text = "Path to VSC certificate pem file" },
    @{ Row = 64; Text = "Path to CA certificate pem file" },
    @{ Row = 65; Text = "XMPP domain used in custom certificates" },
    @{ Row = 66; Text = "Name of the credentials set for the vsc" },
    @{ Row = 68; Text = "Cpuset information for cpu pinning on KVM. For example, VSC requires 4 cores and sample values will be of the form [ 0, 1, 2, 3 ] (List items separated by comma.)" },
    @{ Row = 69; Text = "Enables hardening configuration on VSC [default: True]" }
)

# Step 1: insert a new blank row at row 18 (pushes "VM name" and everything
# below down by one row; cell values/styles/merges/data validations shift
# automatically, but comments stay pinned to their original row, so we fix
# those up afterwards).
$ws.Rows("18:18").Insert()

# Step 2: populate the new row 18 with the "Router ID" label, matching the
# style used by sibling single-column label rows (e.g. row 17 "System IP
# address").
$ws.Range("A18").Value = "Router ID"
$ws.Range("A18").Style = $ws.Range("A17").Style
$ws.Range("B18:C18").Style = $ws.Range("B17:C17").Style

# Step 3: shift every pre-existing comment at/after the old row 18 down by one
# row by rewriting its text in place (applied from the pre-insert snapshot).
foreach ($item in $commentShift) {
    $cell = $ws.Range("A" + $item.Row)
    $cmt = $cell.Comment
    if ($cmt -eq $null) {
        $cmt = $cell.AddComment($item.Text)
    } else {
        $cmt.Text($item.Text)
    }
}

# Step 4: give the new "Router ID" row its own comment describing the field.
$c18 = $ws.Range("A18").Comment
if ($c18 -eq $null) {
    $ws.Range("A18").AddComment("Required IPv4 address when using an IPv6 system IP address [default: (system_ip)]")
} else {
    $c18.Text("Required IPv4 address when using an IPv6 system IP address [default: (system_ip)]")
}
